# Daily attendance processing - 2025-12-09 23:26:44
# Reorders the "Recorded By" (column G) names for rows that were
# recorded by "System" together with a backup/admin account, so that
# the non-"System" account is listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, backup@backdoor.com, system") {
        $cell.Value2 = "system, System, backup@backdoor.com"
    }
}
